$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITERACION1")
Write-Output "noop"
